$wb = $excel.ActiveWorkbook

$oldId   = "645a5244-477f-41e9-8df2-c9bc9ab3ea47"
$newId   = "9e84174e-27f1-440a-a085-90b014160416"
$newHash = "b059b3d731c47b5f52c36c33e6b5ab62a3164707"

$ghBase       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/beaa13ef3e5483f1c1da9f5f50e4d513f01776b9/e2e/$oldId.md"
$ghZhBase     = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5bd4fa5d42dc047217cdb38d021b220ee3509b2f/e2e/$oldId.md"
$ghDeBase     = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5d86c42e5dc66d80d862d18f4557db40b716142/e2e/$oldId.md"

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("G2").Value = "2016-09-06 23:14:15"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ghBase, "", "", "e2e\$newId.md")

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-06 23:14:09"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $ghBase, "", "", "$newId.md")

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-06 23:14:15"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $ghBase, "", "", "$newId.md")
